$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 7 with uncertainty fluid-fill percentage data
$ws.Range("A7").Value = "UNC_Fluid_Fill"
$ws.Range("B7").Value = 0.8
$ws.Range("C7").Value = "% fluid fill"

# Register the new defined name pointing at the new cell
$wb.Names.Add("UNC_Fluid_Fill", "=Sheet1!`$B`$7")

# Update the fluidMass formula to incorporate the uncertainty fill factor
$ws.Range("B6").Formula = "=(container_height*3.14159*(container_diameter/2)^2)/1000*UNC_Fluid_Fill"

# Widen column A to fit the new label
$ws.Columns("A").ColumnWidth = 16.28515625

# Move the active selection to the newly added cell
$ws.Range("B7").Select()
